$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column header for the new "Level" column
$ws.Range("C1").Value2 = "Level"

# Per-recipient level classification (rows 2-101, matching A2:A101/B2:B101)
$levels = @(
    "Federal",
    "State",
    "Regional",
    "Federal",
    "Local",
    "Local",
    "Local",
    "Local",
    "Local",
    "Local",
    "Local",
    "State",
    "Regional",
    "Local",
    "Local",
    "Regional",
    "Regional",
    "Tribal",
    "Local",
    "County",
    "Regional",
    "Regional",
    "Local",
    "Local",
    "Local",
    "Local",
    "Regional",
    "Local",
    "Local",
    "Tribal",
    "Tribal",
    "County",
    "County",
    "Federal",
    "Local",
    "Local",
    "Tribal",
    "Regional",
    "Local",
    "Local",
    "Regional",
    "Local",
    "Local",
    "County",
    "Tribal",
    "Local",
    "State",
    "Local",
    "Local",
    "Local",
    "Local",
    "Local",
    "Local",
    "Local",
    "Local",
    "Tribal",
    "Local",
    "Tribal",
    "Local",
    "Local",
    "Local",
    "Local",
    "Local",
    "State",
    "Local",
    "Local",
    "Local",
    "Local",
    "State",
    "Local",
    "Local",
    "Local",
    "Local",
    "Local",
    "Local",
    "Local",
    "Local",
    "Local",
    "Local",
    "Local",
    "Local",
    "Local",
    "Local",
    "Local",
    "Local",
    "Local",
    "Local",
    "Local",
    "Local",
    "Local",
    "Local",
    "County",
    "Regional",
    "Local",
    "Local",
    "Local",
    "Local",
    "Local",
    "Local",
    "State"
)

for ($i = 0; $i -lt $levels.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value2 = $levels[$i]
}

# Match author's formatting/view tweaks from the same edit
$ws.Columns.Item(2).ColumnWidth = 33.6640625
$ws.Range("C102").Select()
